# feat: add shipping address
#
# - Address sheet: insert a new "company_name" column after last_name,
#   fix the "fulll_state" -> "full_state" header typo, add sample data.
# - User sheet: add a new user row's company values (ZGAutomation / testfifteen).
# - Selection / active-tab bookkeeping: Address becomes the active sheet/tab.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Address sheet: insert "company_name" column (D) right after last_name (C)
# ---------------------------------------------------------------------------
$wsAddr = $wb.Worksheets.Item("Address")
$wsAddr.Columns("D").Insert()
$wsAddr.Range("D1").Value = "company_name"

# ---------------------------------------------------------------------------
# User sheet: fill in the new user's company / test values
# ---------------------------------------------------------------------------
$wsUser = $wb.Worksheets.Item("User")
$wsUser.Range("B7").Value = "ZGAutomation"
$wsUser.Range("C7").Value = "testfifteen"

# ---------------------------------------------------------------------------
# Address sheet: fix header typo (fulll_state -> full_state); this column
# shifted from H to I after the company_name column was inserted.
# ---------------------------------------------------------------------------
$wsAddr.Range("I1").Value = "full_state"

# ---------------------------------------------------------------------------
# Best-effort cosmetic column widths (close to the authored workbook)
# ---------------------------------------------------------------------------
$wsUser.Columns("B").ColumnWidth = 13.25
$wsUser.Columns("C").ColumnWidth = 9.583333333333334
$wsAddr.Columns("D").ColumnWidth = 14.75

# ---------------------------------------------------------------------------
# Selection / active sheet bookkeeping
# ---------------------------------------------------------------------------
$wsUser.Range("D25").Select()

$wsAddr.Activate()
$wsAddr.Range("I1").Select()

Write-Host "Applied shipping-address edits"
